$d = $word.ActiveDocument

function Split-ParagraphIntoWordRuns($paragraphIndex, $styleName, $expectedText, $words) {
    $p = $d.Paragraphs($paragraphIndex)
    $range = $p.Range
    $style = $range.ParagraphFormat.Style.NameLocal
    # Paragraph.Range includes the trailing paragraph mark; strip it off so
    # the replacement XML only covers the paragraph mark + existing content,
    # which InsertXML needs (it expects a full <w:p> element to replace the
    # paragraph, paragraph mark included).
    if ($style -ne $styleName) {
        throw "Paragraph $paragraphIndex style mismatch: expected $styleName but found $style"
    }
    if ($range.Text -ne ($expectedText + [char]13)) {
        throw "Paragraph $paragraphIndex text mismatch: found [$($range.Text)]"
    }

    $runsXml = ""
    for ($i = 0; $i -lt $words.Count; $i++) {
        if ($i -gt 0) {
            $runsXml += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
        }
        $escaped = $words[$i] -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
        $runsXml += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }

    $xml = '<w:p><w:pPr><w:pStyle w:val="' + $styleName + '"/></w:pPr>' + $runsXml + '</w:p>'
    $range.InsertXML($xml)
}

Split-ParagraphIntoWordRuns 1 "Title" "Answers: Introduction to complex numbers" `
    @("Answers:", "Introduction", "to", "complex", "numbers")

Split-ParagraphIntoWordRuns 2 "Author" "Tom Coleman" `
    @("Tom", "Coleman")

Split-ParagraphIntoWordRuns 4 "Abstract" "Answers to questions relating to the guide on introduction to complex numbers." `
    @("Answers", "to", "questions", "relating", "to", "the", "guide", "on", "introduction", "to", "complex", "numbers.")

Write-Host "Done"
